$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1 ("Title" slide) - shape 6 holds the stacked title-block text.
# ---------------------------------------------------------------------------
$s1  = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(6)
$tr1 = $sh1.TextFrame.TextRange

# Paragraph "University of Edinburgh": merge the separate "of " / "Edinburgh"
# runs into a single "of Edinburgh" run (same visible text, fewer runs).
$paraUni = $tr1.Paragraphs(4)
$ofEdinburgh = $paraUni.Characters(12, 12)
$ofEdinburgh.Text = "of Edinburgh"

# Paragraph "Ó  Diarmuid McDonnell, Vernon Gayle": merge the " " and
# "Diarmuid " runs into a single " Diarmuid " run.
$paraAuthor = $tr1.Paragraphs(9)
$diarmuid = $paraAuthor.Characters(3, 10)
$diarmuid.Text = " Diarmuid "

# ---------------------------------------------------------------------------
# Slide 5 ("How does it work?") - shape 2 is the bulleted content body.
# ---------------------------------------------------------------------------
$s5  = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(2)
$tr5 = $sh5.TextFrame.TextRange

# First bullet: merge "Read the brief " / "in the 'hackathon' folder " /
# "on the workshop " runs into one run, leaving "Github" and
# " repository." as their own runs.
$paraRead = $tr5.Paragraphs(1)
$readIntro = $paraRead.Characters(1, 57)
$readIntro.Text = "Read the brief in the " + [char]0x2018 + "hackathon" + [char]0x2019 + " folder on the workshop "

# Third bullet: "...there isn't a solution, so..." becomes
# "...there isn't an official solution, so..." with the new wording its
# own run, sandwiched between the unchanged leading/trailing runs.
$paraStuck = $tr5.Paragraphs(3)
$solution = $paraStuck.Characters(92, 10)
$solution.Text = "an official solution"
